$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 33 (19 marras) ---
# A second time block was added, so "Kello" now lists two ranges and needs
# wrapping like the other multi-range time cells.
$ws.Range("B33").WrapText = $true
$ws.Range("B33").Value = "16.30-17.30, 18.00-19.00"

# --- Row 32 (18 marras) ---
# "Huomiot koodista" note gets extended with two extra sentences.
$ws.Range("E32").Value = "Pikku siistimistä. Haluan katsoa tämän collision pyräyksen innoittamana nyt sitten onnistuisiko vaihtaa esimerkiksi siihen kasipuuksi alla oleva tietorakenne, ja miten iso aikapanostusta vaatisi. Ei ehkä kurssin aikana, ja muutenkin tässä tulee kyllä uutta vielä aika reippaasti tulevina tunteina. Katsotaan mihin riittää aika ja energia."

# The "Oppimisen sisältö" text is extended with the new collision-detection work.
$ws.Range("C33").Value = "Rajaavasta rakenteesta poistaminen, nopea johdanto partiotiointirakenteisiin ja Törmäystarkastelun ABC:tä millingtonin kanssa. Pääsin kahden pallon törmäykseen, s 249-279."

# Hours logged for the newly-added second time block.
$ws.Range("G33").Value = 2

# Row heights grow to fit the longer wrapped text.
$ws.Rows.Item(32).RowHeight = 145
$ws.Rows.Item(33).RowHeight = 101.5

# --- View state ---
$excel.ActiveWindow.ScrollRow = 26
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D33").Select()
